$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph. Rebuild the title paragraph + the new
#    paragraph together via InsertXML so the new paragraph gets the
#    exact run layout: an empty leading run, a bold "Meta description"
#    run, and a plain run with the rest of the sentence.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Bomber Squad Free: Unique 5-Reel Slot with Progressive Jackpot</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read a review of Bomber Squad, a five-reel online slot with 1024 ways to win, a thrilling storyline, and a unique bank robbery theme. Play for free now!</w:t></w:r></w:p>'

$titleRange.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph that used to live near
#    the end of the document (right before the closing italic blurb).
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Play Bomber Squad Free: Unique 5-Reel Slot with Progressive Jackpot" -and $i -gt 1) {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 3) Swap out the final italic paragraph's text for the new image
#    prompt copy. Scope the Find/Replace to just that last paragraph
#    so the identical sentence that now also lives in the "Meta
#    description" paragraph (inserted in step 1) is left untouched.
# ------------------------------------------------------------------
$oldText = "Read a review of Bomber Squad, a five-reel online slot with 1024 ways to win, a thrilling storyline, and a unique bank robbery theme. Play for free now!"
$newText = "Create a feature image that captures the excitement and adventure of Bomber Squad by SimplePlay. The image should be in a cartoon style and should feature a happy Maya warrior wearing glasses. The Maya warrior should be holding a dynamite stick and standing in front of a bank vault door while gangsters shoot their machine guns in the background. Use bright colors and bold outlines to make the image stand out and convey the high-energy and thrilling experience of playing this slot game. The image should make players curious and excited to try out Bomber Squad."

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
